# Update "想去人数" (want-to-go count) values in the 展览 (Exhibition) sheet
# and the corresponding rows in the 全部类型 (All Types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for the 展览 sheet (column F)
$exhibitionUpdates = @{
    3  = 163
    4  = 7910
    5  = 95
    10 = 460
    13 = 447
    15 = 72
    16 = 28
    17 = 5799
    18 = 173
    19 = 250
    20 = 1642
    21 = 231
    22 = 368
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value updates for the 全部类型 sheet (column F)
$allTypesUpdates = @{
    3  = 163
    4  = 7910
    5  = 95
    10 = 460
    13 = 447
    15 = 72
    16 = 28
    18 = 5799
    20 = 173
    21 = 250
    22 = 1642
    23 = 231
    24 = 368
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
